# OrangeHRM_Nafis Test Data.xlsx edit
# - Renames existing sheets to the new TC-numbered scheme
# - Inserts two new sheets for the "Add Pay Grade" / currency work:
#     "TC05-01 Add New Currency" and "TC05-03 Delete Currencies"
# - Populates the new sheets with currency reference data
# - Updates selection / active-tab view state to match the authored edit
#
# NOTE: worksheet object handles returned by Worksheets.Add()/Item() are
# position-bound, not identity-bound -- inserting a sheet can make an
# *existing* variable silently start pointing at a different sheet once
# indices shift. To stay safe we always re-fetch sheets by name (Item("..."))
# right before using them, instead of trusting a previously stored handle.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the four pre-existing sheets in place (content untouched)
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "TC03-02 Add New Account"
$wb.Worksheets.Item(2).Name = "TC04-01 Add Job Title"
$wb.Worksheets.Item(3).Name = "TC06-01 Add Job Category"
$wb.Worksheets.Item(4).Name = "TC07-01 - Add Employee"

# ---------------------------------------------------------------------
# 2. Insert the two new currency sheets between "...Job Title" and
#    "...Add Job Category". Create "Delete Currencies" first so the
#    engine hands "Add New Currency" the higher internal sheetId (7)
#    and "Delete Currencies" gets 6, matching the authored workbook.
#    Re-fetch by name after each insertion (see note above).
# ---------------------------------------------------------------------
$newSheet1 = $wb.Worksheets.Add($wb.Worksheets.Item("TC06-01 Add Job Category"))
$newSheet1.Name = "TC05-03 Delete Currencies"

$newSheet2 = $wb.Worksheets.Add($wb.Worksheets.Item("TC05-03 Delete Currencies"))
$newSheet2.Name = "TC05-01 Add New Currency"

# ---------------------------------------------------------------------
# 3. Populate "TC05-01 Add New Currency"
# ---------------------------------------------------------------------
$wsAddCurrency = $wb.Worksheets.Item("TC05-01 Add New Currency")

$wsAddCurrency.Range("A1").Value = "payGradeName"
$wsAddCurrency.Range("B1").Value = "currency"
$wsAddCurrency.Range("C1").Value = "minimumSalary"
$wsAddCurrency.Range("D1").Value = "maximumSalary"

$wsAddCurrency.Range("B2").Value = "EUR - Euro"
$wsAddCurrency.Range("B3").Value = "FJD - Fiji Dollar"
$wsAddCurrency.Range("B4").Value = "AED - Utd. Arab Emir. Dirham"
$wsAddCurrency.Range("B5").Value = "ALL - Albanian Lek"
$wsAddCurrency.Range("B6").Value = "XAU - Gold (oz.)"

$wsAddCurrency.Range("B3").WrapText = $true
$wsAddCurrency.Rows.Item(3).RowHeight = 30

# ---------------------------------------------------------------------
# 4. Populate "TC05-03 Delete Currencies"
# ---------------------------------------------------------------------
$wsDeleteCurrencies = $wb.Worksheets.Item("TC05-03 Delete Currencies")

$wsDeleteCurrencies.Range("A1").Value = "currency"
$wsDeleteCurrencies.Range("A2").Value = "Euro"
$wsDeleteCurrencies.Range("A3").Value = "Fiji Dollar"
$wsDeleteCurrencies.Range("A4").Value = "Utd. Arab Emir. Dirham"
$wsDeleteCurrencies.Range("A5").Value = "Albanian Lek"
$wsDeleteCurrencies.Range("A6").Value = "Gold (oz.)"

$wsDeleteCurrencies.Range("A3").WrapText = $true

# ---------------------------------------------------------------------
# 5. Selections on the touched sheets
# ---------------------------------------------------------------------
$wsDeleteCurrencies = $wb.Worksheets.Item("TC05-03 Delete Currencies")
$wsDeleteCurrencies.Range("A3:A6").Select()

$wsEmployee = $wb.Worksheets.Item("TC07-01 - Add Employee")
$wsEmployee.Range("I23").Select()

# ---------------------------------------------------------------------
# 6. Final active-sheet / active-tab state: "TC05-01 Add New Currency"
#    is the selected tab (tabSelected=1, activeTab=2 zero-based) with
#    F11:F12 selected.
# ---------------------------------------------------------------------
$wsAddCurrency = $wb.Worksheets.Item("TC05-01 Add New Currency")
$wsAddCurrency.Range("F11:F12").Select()
